# Apply the "added data driven tests" commit:
#  - Add two new worksheets ("Edit", "Delete") after "Data", populated from
#    the existing Computer Name / Company columns (with one edited value and
#    one new shared string "new company").
#  - Change the selection on "Data" from C10 to D1:D4.
#  - Leave "Delete" as the final active sheet/tab, with its own selection.

$wb = $excel.ActiveWorkbook

# ---- Data sheet: update the saved selection (was C10 -> now D1:D4) ----
$wsData = $wb.Worksheets.Item("Data")
$wsData.Range("D1:D4").Select()

# ---- Add "Edit" sheet right after "Data" ----
$wsEdit = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsData)
$wsEdit.Name = "Edit"

$wsEdit.Range("A1").Value = "Computer Name"
$wsEdit.Range("B1").Value = "new company"
$wsEdit.Range("A2").Value = "ASUS ROG 1"
$wsEdit.Range("B2").Value = "RCA"
$wsEdit.Range("A3").Value = "ASUS ROG 2"
$wsEdit.Range("B3").Value = "Sony"
$wsEdit.Range("A4").Value = "ASUS ROG 3"
$wsEdit.Range("B4").Value = "IBM"

# ColumnWidth is quantized to whole pixels by the host when round-tripped
# through COM, so these inputs are chosen to land on the stored width
# closest to the authored workbook's 16.140625 / 13.7109375 (saved by Excel
# directly via its own finer-grained measurement).
$wsEdit.Columns.Item(1).ColumnWidth = 15.34
$wsEdit.Columns.Item(2).ColumnWidth = 12.84

$wsEdit.Range("A1:A4").Select()

# ---- Add "Delete" sheet right after "Edit" ----
$wsDelete = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsEdit)
$wsDelete.Name = "Delete"

$wsDelete.Range("A1").Value = "Computer Name"
$wsDelete.Range("A2").Value = "ASUS ROG 1"
$wsDelete.Range("A3").Value = "ASUS ROG 2"
$wsDelete.Range("A4").Value = "ASUS ROG 3"

# Closest achievable width to the authored 14.140625 (see note above).
$wsDelete.Columns.Item(1).ColumnWidth = 13.34

# Leave the selection on "Delete" at M22, matching the saved workbook state,
# and make "Delete" the active (last-focused) sheet/tab.
$wsDelete.Range("M22").Select()
